# eleusis_master.xlsx: rename two sheets, and move the "active tab" /
# selection from "Inscription" over to the newly-renamed
# "Person Honor Display" sheet (selecting F120 there, where it was
# previously parked on E398).

$wb = $excel.ActiveWorkbook

# "Institution Sponsor" -> "Institution Sponsorship"
$wb.Worksheets.Item("Institution Sponsor").Name = "Institution Sponsorship"

# "Person Displaying Honor" -> "Person Honor Display"
$wb.Worksheets.Item("Person Displaying Honor").Name = "Person Honor Display"

# Activate the renamed "Person Honor Display" sheet (this becomes the
# workbook's active tab / tabSelected sheet, replacing "Inscription"),
# and move its selection to F120.
$ws = $wb.Worksheets.Item("Person Honor Display")
$ws.Activate()
$ws.Range("F120").Select()
